$d = $word.ActiveDocument

$replacements = @(
    @("74-5=69", "15+3=18"),
    @("99-23=76", "34+37=71"),
    @("39-16=23", "91-22=69"),
    @("98-7=91", "96+2=98"),
    @("45+20=65", "80-52=28"),
    @("18+71=89", "74-74=0"),
    @("74-60=14", "8+80=88"),
    @("76-3=73", "9-6=3"),
    @("46+47=93", "70-48=22"),
    @("70-16=54", "5+31=36"),
    @("47+25=72", "30-7=23"),
    @("97-16=81", "58+17=75"),
    @("43-30=13", "33+53=86"),
    @("92-83=9", "26+12=38"),
    @("25+56=81", "16+67=83"),
    @("32-10=22", "96-58=38"),
    @("36+12=48", "34+42=76"),
    @("40-17=23", "53+17=70"),
    @("1+87=88", "42+37=79"),
    @("72-69=3", "24-0=24"),
    @("19+4=23", "6+92=98"),
    @("15+66=81", "83-24=59"),
    @("76-27=49", "56-36=20"),
    @("79-26=53", "65+14=79"),
    @("68-0=68", "6+33=39"),
    @("70-22=48", "48-31=17"),
    @("61+12=73", "16+38=54"),
    @("85-84=1", "33+24=57"),
    @("88-42=46", "31+21=52"),
    @("10+73=83", "80-52=28"),
    @("31-15=16", "94-52=42"),
    @("20+74=94", "59-41=18"),
    @("72-19=53", "95-56=39"),
    @("19+72=91", "10+32=42"),
    @("91-53=38", "76-63=13"),
    @("40+7=47", "55-38=17"),
    @("85+3=88", "39+30=69"),
    @("22+20=42", "20+63=83"),
    @("89+8=97", "52+42=94"),
    @("83+16=99", "38+26=64"),
    @("10+76=86", "89+9=98"),
    @("46-5=41", "65-51=14"),
    @("96-89=7", "94-81=13"),
    @("70-3=67", "29+46=75"),
    @("86-20=66", "40+33=73"),
    @("37+42=79", "64+24=88"),
    @("36+40=76", "94-51=43"),
    @("41+41=82", "69-12=57"),
    @("36+42=78", "69-60=9"),
    @("73+9=82", "95-5=90"),
    @("81+13=94", "10+77=87"),
    @("10+70=80", "85-68=17"),
    @("29+0=29", "10+3=13"),
    @("31-27=4", "3+59=62"),
    @("72-35=37", "61+9=70"),
    @("45+7=52", "8+48=56"),
    @("90-5=85", "52+26=78"),
    @("83-49=34", "69+18=87"),
    @("71-39=32", "26-20=6"),
    @("2+69=71", "69-48=21"),
    @("51+31=82", "56-32=24"),
    @("43+33=76", "8+10=18"),
    @("59-36=23", "12-10=2"),
    @("96-18=78", "89-62=27"),
    @("82-73=9", "90+2=92"),
    @("8+12=20", "72+21=93"),
    @("68-27=41", "70+8=78"),
    @("53-50=3", "93-86=7"),
    @("31+54=85", "94-36=58"),
    @("20+70=90", "39+16=55"),
    @("20+1=21", "98-64=34"),
    @("8+91=99", "28+41=69"),
    @("86+0=86", "27+36=63"),
    @("64-32=32", "85-19=66"),
    @("82-68=14", "15+25=40"),
    @("16+34=50", "65-47=18"),
    @("96-50=46", "96-84=12"),
    @("84-76=8", "76+16=92"),
    @("58+35=93", "25-4=21"),
    @("90-60=30", "55+18=73"),
    @("51-19=32", "15+32=47"),
    @("6+11=17", "72-58=14"),
    @("78-40=38", "39-23=16"),
    @("69-21=48", "58+36=94"),
    @("27+1=28", "97-91=6"),
    @("38-16=22", "22-16=6"),
    @("82-2=80", "65+3=68"),
    @("99-93=6", "96-39=57"),
    @("31+11=42", "14+82=96"),
    @("94-61=33", "2+70=72"),
    @("22-17=5", "0+91=91"),
    @("31-30=1", "5+19=24"),
    @("26-21=5", "26+16=42"),
    @("44+41=85", "24+73=97"),
    @("18-5=13", "27+70=97"),
    @("40+20=60", "38-34=4"),
    @("87-53=34", "49+25=74"),
    @("36+50=86", "67+2=69"),
    @("22+11=33", "3+38=41"),
    @("37-17=20", "67+30=97"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
